# SMARTNODES.docx (Portuguese) — content update
#
# Replaces the two "SmartNodes enable the following services" bullet
# paragraphs with their revised copy (new wording + refreshed run/paragraph
# formatting + switched numbering list), and renumbers the "smarthosting"
# bookmark id from 0 to 1 to mirror the upstream resave.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml($paragraph, [string]$innerParagraphXml) {
    # Use the paragraph's full range (including its end-of-paragraph mark) so
    # InsertXML replaces the whole paragraph in place (props + runs) instead
    # of appending a sibling run ahead of the old paragraph mark.
    $rng = $paragraph.Range
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document ' + $wNs + '>' + $innerParagraphXml + '</w:document></pkg:xmlData>' +
           '</pkg:part></pkg:package>'
    $rng.InsertXML($pkg) | Out-Null
}

# --- Locate the two target bullet paragraphs by their current text ---------
$instantPayPara = $null
$moreServicesPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "InstantPay*instant transactions*") {
        $instantPayPara = $i
    } elseif ($t -like "Will have more services added later*") {
        $moreServicesPara = $i
    }
}

# --- Paragraph: InstantPay bullet ------------------------------------------
$instantPayXml = '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
    '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
    '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
    '<w:textAlignment w:val="baseline"/>' +
    '<w:rPr><w:rFonts w:ascii="Open Sans" w:eastAsia="Times New Roman" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
    '<w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/>' +
    '<w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr>' +
    '<w:t xml:space="preserve">InstantPay (Instant Transactions): Allows for SmartCash transactions to be locked in about a second. No risk of double spending a transaction, so the receiver can trust that transaction immediately.</w:t></w:r></w:p>'

Set-ParagraphXml $d.Paragraphs.Item($instantPayPara) $instantPayXml

# --- Paragraph: SmartRewards bullet ----------------------------------------
$smartRewardsXml = '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
    '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
    '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
    '<w:textAlignment w:val="baseline"/>' +
    '<w:rPr><w:rFonts w:ascii="Open Sans" w:eastAsia="Times New Roman" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
    '<w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/>' +
    '<w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr>' +
    '<w:t xml:space="preserve">SmartRewards: SmartRewards are calculated by the SmartNodes to allow for distribution to be handled automatically by the block rewards.</w:t></w:r></w:p>'

Set-ParagraphXml $d.Paragraphs.Item($moreServicesPara) $smartRewardsXml

Write-Output "done"
